$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 and 7 (the "Vostro / Dell / DKS" device entries), shifting
# subsequent rows up.
$ws.Rows("6:7").Delete()

# Page setup as recorded for this worksheet (adds the <pageSetup> element).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the selection shown in the saved view.
$ws.Range("E16").Select()
